$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")
$ws.Range("H2").Value = "Static expenses:"
